$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 982.4
$ws.Range("J127").Value = 1446.2
$ws.Range("L127").Value = 4338.6
$ws.Range("N127").Value = -14258.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3032192
$ws.Range("I137").Value = 3334914.5
$ws.Range("K137").Value = 10004743.5
$ws.Range("M137").Value = -10002193.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16064.529
$ws.Range("I32").Value = 16329.054
$ws.Range("J32").Value = 14080.6
$ws.Range("K32").Value = 16329.054
$ws.Range("L32").Value = 14080.6
$ws.Range("M32").Value = -16042.054
$ws.Range("N32").Value = -14654.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 30364450
$ws.Range("I61").Value = 38500604
$ws.Range("J61").Value = 144447.28
$ws.Range("K61").Value = 38500604
$ws.Range("L61").Value = 144447.28
$ws.Range("M61").Value = -38500392
$ws.Range("N61").Value = -144871.28

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8402006
$ws.Range("I74").Value = 12551579
$ws.Range("J74").Value = 102861
$ws.Range("K74").Value = 12551579
$ws.Range("L74").Value = 102861
$ws.Range("M74").Value = -12550705
$ws.Range("N74").Value = -104609

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 8402006
$ws.Range("I77").Value = 12551579
$ws.Range("J77").Value = 102861
$ws.Range("K77").Value = 62757895
$ws.Range("L77").Value = 514305
$ws.Range("M77").Value = -62753527
$ws.Range("N77").Value = -523041

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3290282.2
$ws.Range("I97").Value = 3906900
$ws.Range("J97").Value = 1653.6666
$ws.Range("K97").Value = 3906900
$ws.Range("L97").Value = 1653.6666
$ws.Range("M97").Value = -3906404
$ws.Range("N97").Value = -2645.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3269840.2
$ws.Range("I122").Value = 1853.3871
$ws.Range("J122").Value = 37039036
$ws.Range("K122").Value = 5560.1613
$ws.Range("L122").Value = 111117108
$ws.Range("M122").Value = -3110.1613
$ws.Range("N122").Value = -111122008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 52454.27
$ws.Range("I132").Value = 37949.48
$ws.Range("J132").Value = 80427.78999999999
$ws.Range("K132").Value = 113848.44
$ws.Range("L132").Value = 241283.37
$ws.Range("M132").Value = -111318.44
$ws.Range("N132").Value = -246343.37

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 30364450
$ws.Range("I136").Value = 38500604
$ws.Range("J136").Value = 144447.28
$ws.Range("K136").Value = 115501812
$ws.Range("L136").Value = 433341.84
$ws.Range("M136").Value = -115499262
$ws.Range("N136").Value = -438441.84

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 240
$ws.Range("I22").Value = 240
$ws.Range("K22").Value = 240
$ws.Range("M22").Value = -67

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 970.2222
$ws.Range("I64").Value = 500
$ws.Range("J64").Value = 1205.3334
$ws.Range("K64").Value = 500
$ws.Range("L64").Value = 1205.3334
$ws.Range("M64").Value = -275
$ws.Range("N64").Value = -1655.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 970.2222
$ws.Range("I67").Value = 500
$ws.Range("J67").Value = 1205.3334
$ws.Range("K67").Value = 500
$ws.Range("L67").Value = 1205.3334
$ws.Range("M67").Value = 280
$ws.Range("N67").Value = -2765.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1042.7931
$ws.Range("I99").Value = 1024.5454
$ws.Range("K99").Value = 1024.5454
$ws.Range("M99").Value = 473.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3119.8
$ws.Range("I134").Value = 2460.8684
$ws.Range("J134").Value = 5206.4165
$ws.Range("K134").Value = 7382.6052
$ws.Range("L134").Value = 15619.2495
$ws.Range("M134").Value = -4847.6052
$ws.Range("N134").Value = -20689.2495

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2571.425
$ws.Range("I31").Value = 1294.9678
$ws.Range("J31").Value = 6968.1113
$ws.Range("K31").Value = 1294.9678
$ws.Range("L31").Value = 6968.1113
$ws.Range("M31").Value = -999.9677999999999
$ws.Range("N31").Value = -7558.1113

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2571.425
$ws.Range("I34").Value = 1294.9678
$ws.Range("J34").Value = 6968.1113
$ws.Range("K34").Value = 1294.9678
$ws.Range("L34").Value = 6968.1113
$ws.Range("M34").Value = -1092.9678
$ws.Range("N34").Value = -7372.1113

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4639.6924
$ws.Range("J94").Value = 1310.3334
$ws.Range("L94").Value = 1310.3334
$ws.Range("N94").Value = -2212.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 829.6667
$ws.Range("I105").Value = 794.93335
$ws.Range("J105").Value = 1003.3333
$ws.Range("K105").Value = 794.93335
$ws.Range("L105").Value = 1003.3333
$ws.Range("M105").Value = 952.06665
$ws.Range("N105").Value = -4497.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 45723.39
$ws.Range("I132").Value = 1982
$ws.Range("K132").Value = 5946
$ws.Range("M132").Value = -3416

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 885.3570999999999
$ws.Range("I122").Value = 285.33334
$ws.Range("J122").Value = 1049
$ws.Range("K122").Value = 2568.00006
$ws.Range("L122").Value = 9441
$ws.Range("M122").Value = -118.0000600000003
$ws.Range("N122").Value = -14341

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 2723.4167
$ws.Range("I123").Value = 1985
$ws.Range("J123").Value = 2969.5557
$ws.Range("K123").Value = 5955
$ws.Range("L123").Value = 8908.667099999999
$ws.Range("M123").Value = -3505
$ws.Range("N123").Value = -13808.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1173.9207
$ws.Range("I131").Value = 508.16666
$ws.Range("J131").Value = 1244
$ws.Range("K131").Value = 1524.49998
$ws.Range("L131").Value = 3732
$ws.Range("M131").Value = 3515.50002
$ws.Range("N131").Value = -13812

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1007.5
$ws.Range("I102").Value = 832.75
$ws.Range("J102").Value = 1240.5
$ws.Range("K102").Value = 832.75
$ws.Range("L102").Value = 1240.5
$ws.Range("M102").Value = 789.25
$ws.Range("N102").Value = -4484.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2330.5898
$ws.Range("I122").Value = 1944.0344
$ws.Range("J122").Value = 3451.6
$ws.Range("K122").Value = 5832.1032
$ws.Range("L122").Value = 10354.8
$ws.Range("M122").Value = -3382.1032
$ws.Range("N122").Value = -15254.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6900
$ws.Range("I126").Value = 6900
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 20700
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -18230
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 120222.82
$ws.Range("I132").Value = 101679
$ws.Range("K132").Value = 305037
$ws.Range("M132").Value = -302507

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 14136
$ws.Range("J141").Value = 21300
$ws.Range("L141").Value = 21300
$ws.Range("N141").Value = -31660

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 37500
$ws.Range("J75").Value = 37500
$ws.Range("L75").Value = 37500
$ws.Range("N75").Value = -39372

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H78").Value = 37500
$ws.Range("J78").Value = 37500
$ws.Range("L78").Value = 112500
$ws.Range("N78").Value = -121860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1990
$ws.Range("I100").Value = 1950
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1950
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1409
$ws.Range("N100").Value = -3082

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 40068.48
$ws.Range("I132").Value = 2412.9473
$ws.Range("J132").Value = 129500.375
$ws.Range("K132").Value = 7238.841899999999
$ws.Range("L132").Value = 388501.125
$ws.Range("M132").Value = -4708.841899999999
$ws.Range("N132").Value = -393561.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 49928.285
$ws.Range("I136").Value = 28215.893
$ws.Range("J136").Value = 210600
$ws.Range("K136").Value = 84647.679
$ws.Range("L136").Value = 631800
$ws.Range("M136").Value = -82097.679
$ws.Range("N136").Value = -636900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 15000
$ws.Range("I96").Value = 11000
$ws.Range("K96").Value = 11000
$ws.Range("M96").Value = -9627

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 46671.953
$ws.Range("I132").Value = 32084
$ws.Range("J132").Value = 85573.164
$ws.Range("K132").Value = 96252
$ws.Range("L132").Value = 256719.492
$ws.Range("M132").Value = -93722
$ws.Range("N132").Value = -261779.492
